$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F1 header "REMARKS" removed (cell kept, just emptied) ---
$ws.Range("F1").ClearContents()

# --- SOLVED? column: flip NO -> YES (and Bad -> Good style) for rows 2,3,5 ---
$ws.Range("E2").Value = "YES"
$ws.Range("E2").Style = "Good"

$ws.Range("E3").Value = "YES"
$ws.Range("E3").Style = "Good"

$ws.Range("E5").Value = "YES"
$ws.Range("E5").Style = "Good"

# --- New row 6: second bot entry ---
$ws.Range("A6").Value = "Bot/Player not totally dead after death."
$ws.Range("B6").Value = "PR0Pawn (presumably)"
$ws.Range("C6").Value = "Genio"
$ws.Range("E6").Value = "NO"
$ws.Range("E6").Style = "Bad"
$ws.Range("F6").Value = "Maybe because Melkar override the deadanimation function?"

# --- column F width / best fit ---
$ws.Columns("F").ColumnWidth = 43.6640625

# --- view state ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F4").Select()

# --- page setup ---
$ws.PageSetup.Orientation = 1
